$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 9-45: Fecha (D), Volumen (J), Precio minimo/maximo/promedio (K/L/M), Origen (O), Precio $/Kg (P)
$ws.Cells.Item(9, 4).Value = 44425
$ws.Cells.Item(9, 10).Value = 200
$ws.Cells.Item(9, 11).Value = 9000
$ws.Cells.Item(9, 12).Value = 9000
$ws.Cells.Item(9, 13).Value = 9000
$ws.Cells.Item(9, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(9, 16).Value = 150

$ws.Cells.Item(10, 4).Value = 44249
$ws.Cells.Item(10, 10).Value = 200
$ws.Cells.Item(10, 11).Value = 8000
$ws.Cells.Item(10, 12).Value = 8000
$ws.Cells.Item(10, 13).Value = 8000
$ws.Cells.Item(10, 15).Value = "Región del Maule"
$ws.Cells.Item(10, 16).Value = 133

$ws.Cells.Item(11, 4).Value = 44252
$ws.Cells.Item(11, 10).Value = 200
$ws.Cells.Item(11, 11).Value = 8000
$ws.Cells.Item(11, 12).Value = 8000
$ws.Cells.Item(11, 13).Value = 8000
$ws.Cells.Item(11, 15).Value = "Región del Maule"
$ws.Cells.Item(11, 16).Value = 133

$ws.Cells.Item(12, 4).Value = 44293
$ws.Cells.Item(12, 10).Value = 200
$ws.Cells.Item(12, 11).Value = 8000
$ws.Cells.Item(12, 12).Value = 8000
$ws.Cells.Item(12, 13).Value = 8000
$ws.Cells.Item(12, 15).Value = "Región del Maule"
$ws.Cells.Item(12, 16).Value = 133

$ws.Cells.Item(13, 4).Value = 44229
$ws.Cells.Item(13, 10).Value = 200
$ws.Cells.Item(13, 11).Value = 8000
$ws.Cells.Item(13, 12).Value = 8000
$ws.Cells.Item(13, 13).Value = 8000
$ws.Cells.Item(13, 15).Value = "Región del Maule"
$ws.Cells.Item(13, 16).Value = 133

$ws.Cells.Item(14, 4).Value = 44424
$ws.Cells.Item(14, 10).Value = 300
$ws.Cells.Item(14, 11).Value = 8000
$ws.Cells.Item(14, 12).Value = 8000
$ws.Cells.Item(14, 13).Value = 8000
$ws.Cells.Item(14, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(14, 16).Value = 133

$ws.Cells.Item(15, 4).Value = 44236
$ws.Cells.Item(15, 10).Value = 300
$ws.Cells.Item(15, 11).Value = 7000
$ws.Cells.Item(15, 12).Value = 7000
$ws.Cells.Item(15, 13).Value = 7000
$ws.Cells.Item(15, 15).Value = "Región del Maule"
$ws.Cells.Item(15, 16).Value = 117

$ws.Cells.Item(16, 4).Value = 44279
$ws.Cells.Item(16, 10).Value = 200
$ws.Cells.Item(16, 11).Value = 8000
$ws.Cells.Item(16, 12).Value = 8000
$ws.Cells.Item(16, 13).Value = 8000
$ws.Cells.Item(16, 15).Value = "Región del Maule"
$ws.Cells.Item(16, 16).Value = 133

$ws.Cells.Item(17, 4).Value = 44414
$ws.Cells.Item(17, 10).Value = 300
$ws.Cells.Item(17, 11).Value = 7000
$ws.Cells.Item(17, 12).Value = 7000
$ws.Cells.Item(17, 13).Value = 7000
$ws.Cells.Item(17, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(17, 16).Value = 117

$ws.Cells.Item(18, 4).Value = 44239
$ws.Cells.Item(18, 10).Value = 200
$ws.Cells.Item(18, 11).Value = 8000
$ws.Cells.Item(18, 12).Value = 8000
$ws.Cells.Item(18, 13).Value = 8000
$ws.Cells.Item(18, 15).Value = "Región del Maule"
$ws.Cells.Item(18, 16).Value = 133

$ws.Cells.Item(19, 4).Value = 44280
$ws.Cells.Item(19, 10).Value = 200
$ws.Cells.Item(19, 11).Value = 7000
$ws.Cells.Item(19, 12).Value = 7000
$ws.Cells.Item(19, 13).Value = 7000
$ws.Cells.Item(19, 15).Value = "Región del Maule"
$ws.Cells.Item(19, 16).Value = 117

$ws.Cells.Item(20, 4).Value = 44299
$ws.Cells.Item(20, 10).Value = 200
$ws.Cells.Item(20, 11).Value = 8000
$ws.Cells.Item(20, 12).Value = 8000
$ws.Cells.Item(20, 13).Value = 8000
$ws.Cells.Item(20, 15).Value = "Región del Maule"
$ws.Cells.Item(20, 16).Value = 133

$ws.Cells.Item(21, 4).Value = 44242
$ws.Cells.Item(21, 10).Value = 300
$ws.Cells.Item(21, 11).Value = 8000
$ws.Cells.Item(21, 12).Value = 8000
$ws.Cells.Item(21, 13).Value = 8000
$ws.Cells.Item(21, 15).Value = "Región del Maule"
$ws.Cells.Item(21, 16).Value = 133

$ws.Cells.Item(22, 4).Value = 44258
$ws.Cells.Item(22, 10).Value = 200
$ws.Cells.Item(22, 11).Value = 8000
$ws.Cells.Item(22, 12).Value = 8000
$ws.Cells.Item(22, 13).Value = 8000
$ws.Cells.Item(22, 15).Value = "Región del Maule"
$ws.Cells.Item(22, 16).Value = 133

$ws.Cells.Item(23, 4).Value = 44243
$ws.Cells.Item(23, 10).Value = 300
$ws.Cells.Item(23, 11).Value = 8000
$ws.Cells.Item(23, 12).Value = 8000
$ws.Cells.Item(23, 13).Value = 8000
$ws.Cells.Item(23, 15).Value = "Región del Maule"
$ws.Cells.Item(23, 16).Value = 133

$ws.Cells.Item(24, 4).Value = 44421
$ws.Cells.Item(24, 10).Value = 200
$ws.Cells.Item(24, 11).Value = 8000
$ws.Cells.Item(24, 12).Value = 8000
$ws.Cells.Item(24, 13).Value = 8000
$ws.Cells.Item(24, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(24, 16).Value = 133

$ws.Cells.Item(25, 4).Value = 44237
$ws.Cells.Item(25, 10).Value = 200
$ws.Cells.Item(25, 11).Value = 8000
$ws.Cells.Item(25, 12).Value = 8000
$ws.Cells.Item(25, 13).Value = 8000
$ws.Cells.Item(25, 15).Value = "Región del Maule"
$ws.Cells.Item(25, 16).Value = 133

$ws.Cells.Item(26, 4).Value = 44273
$ws.Cells.Item(26, 10).Value = 200
$ws.Cells.Item(26, 11).Value = 8000
$ws.Cells.Item(26, 12).Value = 8000
$ws.Cells.Item(26, 13).Value = 8000
$ws.Cells.Item(26, 15).Value = "Región del Maule"
$ws.Cells.Item(26, 16).Value = 133

$ws.Cells.Item(27, 4).Value = 44257
$ws.Cells.Item(27, 10).Value = 150
$ws.Cells.Item(27, 11).Value = 8000
$ws.Cells.Item(27, 12).Value = 8000
$ws.Cells.Item(27, 13).Value = 8000
$ws.Cells.Item(27, 15).Value = "Región del Maule"
$ws.Cells.Item(27, 16).Value = 133

$ws.Cells.Item(28, 4).Value = 44295
$ws.Cells.Item(28, 10).Value = 200
$ws.Cells.Item(28, 11).Value = 8000
$ws.Cells.Item(28, 12).Value = 8000
$ws.Cells.Item(28, 13).Value = 8000
$ws.Cells.Item(28, 15).Value = "Región del Maule"
$ws.Cells.Item(28, 16).Value = 133

$ws.Cells.Item(29, 4).Value = 44298
$ws.Cells.Item(29, 10).Value = 200
$ws.Cells.Item(29, 11).Value = 8000
$ws.Cells.Item(29, 12).Value = 8000
$ws.Cells.Item(29, 13).Value = 8000
$ws.Cells.Item(29, 15).Value = "Región del Maule"
$ws.Cells.Item(29, 16).Value = 133

$ws.Cells.Item(30, 4).Value = 44274
$ws.Cells.Item(30, 10).Value = 150
$ws.Cells.Item(30, 11).Value = 7000
$ws.Cells.Item(30, 12).Value = 7000
$ws.Cells.Item(30, 13).Value = 7000
$ws.Cells.Item(30, 15).Value = "Región del Maule"
$ws.Cells.Item(30, 16).Value = 117

$ws.Cells.Item(31, 4).Value = 44231
$ws.Cells.Item(31, 10).Value = 250
$ws.Cells.Item(31, 11).Value = 8000
$ws.Cells.Item(31, 12).Value = 8000
$ws.Cells.Item(31, 13).Value = 8000
$ws.Cells.Item(31, 15).Value = "Región del Maule"
$ws.Cells.Item(31, 16).Value = 133

$ws.Cells.Item(32, 4).Value = 44278
$ws.Cells.Item(32, 10).Value = 200
$ws.Cells.Item(32, 11).Value = 8000
$ws.Cells.Item(32, 12).Value = 8000
$ws.Cells.Item(32, 13).Value = 8000
$ws.Cells.Item(32, 15).Value = "Región del Maule"
$ws.Cells.Item(32, 16).Value = 133

$ws.Cells.Item(33, 4).Value = 44389
$ws.Cells.Item(33, 10).Value = 300
$ws.Cells.Item(33, 11).Value = 12000
$ws.Cells.Item(33, 12).Value = 12000
$ws.Cells.Item(33, 13).Value = 12000
$ws.Cells.Item(33, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(33, 16).Value = 200

$ws.Cells.Item(34, 4).Value = 44251
$ws.Cells.Item(34, 10).Value = 200
$ws.Cells.Item(34, 11).Value = 8000
$ws.Cells.Item(34, 12).Value = 8000
$ws.Cells.Item(34, 13).Value = 8000
$ws.Cells.Item(34, 15).Value = "Región del Maule"
$ws.Cells.Item(34, 16).Value = 133

$ws.Cells.Item(35, 4).Value = 44250
$ws.Cells.Item(35, 10).Value = 200
$ws.Cells.Item(35, 11).Value = 8000
$ws.Cells.Item(35, 12).Value = 8000
$ws.Cells.Item(35, 13).Value = 8000
$ws.Cells.Item(35, 15).Value = "Región del Maule"
$ws.Cells.Item(35, 16).Value = 133

$ws.Cells.Item(36, 4).Value = 44305
$ws.Cells.Item(36, 10).Value = 200
$ws.Cells.Item(36, 11).Value = 8000
$ws.Cells.Item(36, 12).Value = 8000
$ws.Cells.Item(36, 13).Value = 8000
$ws.Cells.Item(36, 15).Value = "Región del Maule"
$ws.Cells.Item(36, 16).Value = 133

$ws.Cells.Item(37, 4).Value = 44294
$ws.Cells.Item(37, 10).Value = 200
$ws.Cells.Item(37, 11).Value = 9000
$ws.Cells.Item(37, 12).Value = 9000
$ws.Cells.Item(37, 13).Value = 9000
$ws.Cells.Item(37, 15).Value = "Región del Maule"
$ws.Cells.Item(37, 16).Value = 150

$ws.Cells.Item(38, 4).Value = 44417
$ws.Cells.Item(38, 10).Value = 300
$ws.Cells.Item(38, 11).Value = 7000
$ws.Cells.Item(38, 12).Value = 7000
$ws.Cells.Item(38, 13).Value = 7000
$ws.Cells.Item(38, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(38, 16).Value = 117

$ws.Cells.Item(39, 4).Value = 44419
$ws.Cells.Item(39, 10).Value = 200
$ws.Cells.Item(39, 11).Value = 7000
$ws.Cells.Item(39, 12).Value = 7000
$ws.Cells.Item(39, 13).Value = 7000
$ws.Cells.Item(39, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(39, 16).Value = 117

$ws.Cells.Item(40, 4).Value = 44245
$ws.Cells.Item(40, 10).Value = 200
$ws.Cells.Item(40, 11).Value = 8000
$ws.Cells.Item(40, 12).Value = 8000
$ws.Cells.Item(40, 13).Value = 8000
$ws.Cells.Item(40, 15).Value = "Región del Maule"
$ws.Cells.Item(40, 16).Value = 133

$ws.Cells.Item(41, 4).Value = 44265
$ws.Cells.Item(41, 10).Value = 200
$ws.Cells.Item(41, 11).Value = 8000
$ws.Cells.Item(41, 12).Value = 8000
$ws.Cells.Item(41, 13).Value = 8000
$ws.Cells.Item(41, 15).Value = "Región del Maule"
$ws.Cells.Item(41, 16).Value = 133

$ws.Cells.Item(42, 4).Value = 44277
$ws.Cells.Item(42, 10).Value = 200
$ws.Cells.Item(42, 11).Value = 8000
$ws.Cells.Item(42, 12).Value = 8000
$ws.Cells.Item(42, 13).Value = 8000
$ws.Cells.Item(42, 15).Value = "Región del Maule"
$ws.Cells.Item(42, 16).Value = 133

$ws.Cells.Item(43, 4).Value = 44309
$ws.Cells.Item(43, 10).Value = 150
$ws.Cells.Item(43, 11).Value = 8000
$ws.Cells.Item(43, 12).Value = 8000
$ws.Cells.Item(43, 13).Value = 8000
$ws.Cells.Item(43, 15).Value = "Región del Maule"
$ws.Cells.Item(43, 16).Value = 133

$ws.Cells.Item(44, 4).Value = 44253
$ws.Cells.Item(44, 10).Value = 200
$ws.Cells.Item(44, 11).Value = 8000
$ws.Cells.Item(44, 12).Value = 8000
$ws.Cells.Item(44, 13).Value = 8000
$ws.Cells.Item(44, 15).Value = "Región del Maule"
$ws.Cells.Item(44, 16).Value = 133

$ws.Cells.Item(45, 4).Value = 44272
$ws.Cells.Item(45, 10).Value = 200
$ws.Cells.Item(45, 11).Value = 8000
$ws.Cells.Item(45, 12).Value = 8000
$ws.Cells.Item(45, 13).Value = 8000
$ws.Cells.Item(45, 15).Value = "Región del Maule"
$ws.Cells.Item(45, 16).Value = 133

# Append new rows 46-51 (same Mercado/Region/Categoria/etc. as the rest of the block)
$ws.Cells.Item(46, 1).Value = 5
$ws.Cells.Item(46, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(46, 3).Value = "Maule"
$ws.Cells.Item(46, 4).Value = 44230
$ws.Cells.Item(46, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(46, 5).Value = 7
$ws.Cells.Item(46, 6).Value = 100112001
$ws.Cells.Item(46, 7).Value = "Berenjena"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 200
$ws.Cells.Item(46, 11).Value = 9000
$ws.Cells.Item(46, 12).Value = 9000
$ws.Cells.Item(46, 13).Value = 9000
$ws.Cells.Item(46, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(46, 15).Value = "Región del Maule"
$ws.Cells.Item(46, 16).Value = 150
$ws.Cells.Item(46, 17).Value = 60
$ws.Cells.Item(46, 18).Value = "Hortaliza"

$ws.Cells.Item(47, 1).Value = 5
$ws.Cells.Item(47, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(47, 3).Value = "Maule"
$ws.Cells.Item(47, 4).Value = 44232
$ws.Cells.Item(47, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(47, 5).Value = 7
$ws.Cells.Item(47, 6).Value = 100112001
$ws.Cells.Item(47, 7).Value = "Berenjena"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 150
$ws.Cells.Item(47, 11).Value = 8000
$ws.Cells.Item(47, 12).Value = 8000
$ws.Cells.Item(47, 13).Value = 8000
$ws.Cells.Item(47, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(47, 15).Value = "Región del Maule"
$ws.Cells.Item(47, 16).Value = 133
$ws.Cells.Item(47, 17).Value = 60
$ws.Cells.Item(47, 18).Value = "Hortaliza"

$ws.Cells.Item(48, 1).Value = 5
$ws.Cells.Item(48, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(48, 3).Value = "Maule"
$ws.Cells.Item(48, 4).Value = 44270
$ws.Cells.Item(48, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(48, 5).Value = 7
$ws.Cells.Item(48, 6).Value = 100112001
$ws.Cells.Item(48, 7).Value = "Berenjena"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 200
$ws.Cells.Item(48, 11).Value = 8000
$ws.Cells.Item(48, 12).Value = 8000
$ws.Cells.Item(48, 13).Value = 8000
$ws.Cells.Item(48, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(48, 15).Value = "Región del Maule"
$ws.Cells.Item(48, 16).Value = 133
$ws.Cells.Item(48, 17).Value = 60
$ws.Cells.Item(48, 18).Value = "Hortaliza"

$ws.Cells.Item(49, 1).Value = 5
$ws.Cells.Item(49, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(49, 3).Value = "Maule"
$ws.Cells.Item(49, 4).Value = 44244
$ws.Cells.Item(49, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(49, 5).Value = 7
$ws.Cells.Item(49, 6).Value = 100112001
$ws.Cells.Item(49, 7).Value = "Berenjena"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 200
$ws.Cells.Item(49, 11).Value = 8000
$ws.Cells.Item(49, 12).Value = 8000
$ws.Cells.Item(49, 13).Value = 8000
$ws.Cells.Item(49, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(49, 15).Value = "Región del Maule"
$ws.Cells.Item(49, 16).Value = 133
$ws.Cells.Item(49, 17).Value = 60
$ws.Cells.Item(49, 18).Value = "Hortaliza"

$ws.Cells.Item(50, 1).Value = 5
$ws.Cells.Item(50, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(50, 3).Value = "Maule"
$ws.Cells.Item(50, 4).Value = 44284
$ws.Cells.Item(50, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(50, 5).Value = 7
$ws.Cells.Item(50, 6).Value = 100112001
$ws.Cells.Item(50, 7).Value = "Berenjena"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 200
$ws.Cells.Item(50, 11).Value = 8000
$ws.Cells.Item(50, 12).Value = 8000
$ws.Cells.Item(50, 13).Value = 8000
$ws.Cells.Item(50, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(50, 15).Value = "Región del Maule"
$ws.Cells.Item(50, 16).Value = 133
$ws.Cells.Item(50, 17).Value = 60
$ws.Cells.Item(50, 18).Value = "Hortaliza"

$ws.Cells.Item(51, 1).Value = 5
$ws.Cells.Item(51, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(51, 3).Value = "Maule"
$ws.Cells.Item(51, 4).Value = 44418
$ws.Cells.Item(51, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(51, 5).Value = 7
$ws.Cells.Item(51, 6).Value = 100112001
$ws.Cells.Item(51, 7).Value = "Berenjena"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 200
$ws.Cells.Item(51, 11).Value = 8000
$ws.Cells.Item(51, 12).Value = 8000
$ws.Cells.Item(51, 13).Value = 8000
$ws.Cells.Item(51, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(51, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(51, 16).Value = 133
$ws.Cells.Item(51, 17).Value = 60
$ws.Cells.Item(51, 18).Value = "Hortaliza"
